$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update grade values (C2:E4)
$ws.Range("C2").Value = 88
$ws.Range("D2").Value = 89
$ws.Range("E2").Value = 78

$ws.Range("C3").Value = 81
$ws.Range("D3").Value = 56
$ws.Range("E3").Value = 90

$ws.Range("C4").Value = 55
$ws.Range("D4").Value = 23
$ws.Range("E4").Value = 88

# Move selection to E4
$ws.Range("E4").Select()
